# "Generate Report for Handback"
#
# The localization-status workbook tracks handoff/handback state for each
# source file, per target language (zh-cn, de-de) plus an "Overview" roll-up.
# This handback run:
#   1. Flips the Status column (shared by Overview + each language sheet)
#      from "Ready for handoff" to "Handed back: in sync with en-US" for the
#      two real source-file rows (row 2 + row 3 — row 4 is the
#      not-localized .localization-config row and is left alone).
#   2. Populates "Latest Target File" (E) / "Latest Handback File" (F) for
#      those same rows with hyperlinked filenames (re-using the same
#      md / xlf targets already linked from columns A / C — the handback
#      artifact is the same file that was handed off).
#   3. Stamps "Latest Handback DateTime" (G) with the handback timestamp
#      (kept distinct per language sheet, matching each sheet's own handoff
#      clock: zh-cn -> 10:40:08, de-de -> 10:40:25).

function Get-HyperlinkAddress($ws, $cellAddress) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $cellAddress) {
            return $h.Address
        }
    }
    return $null
}

function Add-MatchingHyperlink($ws, $targetCell, $sourceCellAddress, $displayText) {
    $url = Get-HyperlinkAddress $ws $sourceCellAddress
    $ws.Hyperlinks.Add($ws.Range($targetCell), $url, "", "", $displayText)
}

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- 1. Status column: rows 2 & 3 everywhere it appears -------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

# --- 2 & 3. Per-language handback columns (E, F, G) ------------------------

# zh-cn : row 2 (7603ec1e...) and row 3 (9e002f56...)
Add-MatchingHyperlink $zhcn "E2" '$A$2' "7603ec1e-16ed-4c46-a31b-4490e925c4ed.md"
Add-MatchingHyperlink $zhcn "F2" '$C$2' "7603ec1e-16ed-4c46-a31b-4490e925c4ed.97943d0bf264e610d0a88fb28c1abb04b7e83f45.zh-cn.xlf"
$zhcn.Range("G2").Value = "2016-03-09 10:40:08"

Add-MatchingHyperlink $zhcn "E3" '$A$3' "9e002f56-c8fa-478b-b7e0-c10584114ee4.md"
Add-MatchingHyperlink $zhcn "F3" '$C$3' "9e002f56-c8fa-478b-b7e0-c10584114ee4.ecbd3fb1a50bc84d25c04905b2af0e5577684159.zh-cn.xlf"
$zhcn.Range("G3").Value = "2016-03-09 10:40:08"

# de-de : row 2 (7603ec1e...) and row 3 (9e002f56...)
Add-MatchingHyperlink $dede "E2" '$A$2' "7603ec1e-16ed-4c46-a31b-4490e925c4ed.md"
Add-MatchingHyperlink $dede "F2" '$C$2' "7603ec1e-16ed-4c46-a31b-4490e925c4ed.97943d0bf264e610d0a88fb28c1abb04b7e83f45.de-de.xlf"
$dede.Range("G2").Value = "2016-03-09 10:40:25"

Add-MatchingHyperlink $dede "E3" '$A$3' "9e002f56-c8fa-478b-b7e0-c10584114ee4.md"
Add-MatchingHyperlink $dede "F3" '$C$3' "9e002f56-c8fa-478b-b7e0-c10584114ee4.ecbd3fb1a50bc84d25c04905b2af0e5577684159.de-de.xlf"
$dede.Range("G3").Value = "2016-03-09 10:40:25"
